$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (UserTable header row): insert "www" before "username" and shift the
#     remaining headers one column to the right (E3..T3 -> F3..U3) ---
$ws.Cells.Item(3,21).Value = "imageRef"    # U3  (was T3)
$ws.Cells.Item(3,20).Value = "biog"        # T3  (was S3)
$ws.Cells.Item(3,19).Value = "noEvents"    # S3  (was R3)
$ws.Cells.Item(3,18).Value = "noServices"  # R3  (was Q3)
$ws.Cells.Item(3,17).Value = "userType"    # Q3  (was P3)
$ws.Cells.Item(3,16).Value = "County"      # P3  (was O3)
$ws.Cells.Item(3,15).Value = "Town"        # O3  (was N3)
$ws.Cells.Item(3,14).Value = "Address3"    # N3  (was M3)
$ws.Cells.Item(3,13).Value = "Address2"    # M3  (was L3)
$ws.Cells.Item(3,12).Value = "Address1"    # L3  (was K3)
$ws.Cells.Item(3,11).Value = "Address"     # K3  (was J3)
$ws.Cells.Item(3,10).Value = "phNo"        # J3  (was I3)
$ws.Cells.Item(3,9).Value  = "cName"       # I3  (was H3)
$ws.Cells.Item(3,8).Value  = "lname"       # H3  (was G3)
$ws.Cells.Item(3,7).Value  = "fname"       # G3  (was F3)
$ws.Cells.Item(3,6).Value  = "username"    # F3  (was E3)
$ws.Cells.Item(3,5).Value  = "www"         # E3  (new)

# --- Row 7 (EventTable data row): add StartTime/EndTime/startDate/endDate/noAttendees ---
$ws.Cells.Item(7,7).Value  = "StartTime"   # G7
$ws.Cells.Item(7,8).Value  = "EndTime"     # H7
$ws.Cells.Item(7,9).Value  = "startDate"   # I7
$ws.Cells.Item(7,10).Value = "endDate"     # J7
$ws.Cells.Item(7,11).Value = "noAttendees" # K7

# --- Row 12 (EventServiceTable data row): drop venueId/Date/StartTime/EndTime/
#     startDate/endDate/noAttendees, keep only id/eventID/serviceID ---
$ws.Cells.Item(12,9).Value = "serviceID"   # I12 (was "venueId")
$ws.Range("J12:P12").ClearContents()

# --- Column S width change ---
$ws.Columns("S").ColumnWidth = 12.6

# --- Sheet view / selection ---
$ws.Range("V7").Select()
